$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): industry / unit / process / carbon / ced / climate / region
# Bold, thin-box border, centered horizontally, top-aligned vertically.
# ---------------------------------------------------------------------------
$headers = @("industry", "unit", "process", "carbon (kg CO2 eq)", "ced (MJ)", "climate change (kg CO2 eq)", "region")

$headerRange = $ws.Range("A1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous
$headerRange.Borders.Weight = 2            # xlThin

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ---------------------------------------------------------------------------
# Data rows
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 1).Value = "electricity general industry"
$ws.Cells.Item(2, 2).Value = "MJ "
$ws.Cells.Item(2, 3).Value = "Electricity General Industry"
$ws.Cells.Item(2, 4).Value = 0.13564138
$ws.Cells.Item(2, 5).Value = 2.44904
$ws.Cells.Item(2, 6).Value = [double]"3.7820415e-06"
$ws.Cells.Item(2, 7).Value = "Global"

$ws.Cells.Item(3, 1).Value = "electricity general industry"
$ws.Cells.Item(3, 2).Value = "MJ "
$ws.Cells.Item(3, 3).Value = "Electricity General domestic use Low Voltage"
$ws.Cells.Item(3, 4).Value = 0.1427804
$ws.Cells.Item(3, 5).Value = 2.5779369
$ws.Cells.Item(3, 6).Value = [double]"3.9810963e-06"
$ws.Cells.Item(3, 7).Value = "Global"

# ---------------------------------------------------------------------------
# Comments on header row describing the data type of each column
# ---------------------------------------------------------------------------
$excel.UserName = "Data Processor"

$ws.Range("A1").AddComment("Data type: Categorical (text)") | Out-Null
$ws.Range("B1").AddComment("Data type: Various (e.g. kg, kWh)") | Out-Null
$ws.Range("C1").AddComment("Data type: Categorical (text)") | Out-Null
$ws.Range("D1").AddComment("Data type: Carbon footprint") | Out-Null
$ws.Range("E1").AddComment("Data type: Cumulative energy demand") | Out-Null
$ws.Range("F1").AddComment("Data type: Climate change impact") | Out-Null
$ws.Range("G1").AddComment("Data type: Categorical (text)") | Out-Null
